$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "Gleison Silva Freire - 1203037" -> "Gleison Freire - 1203037"
#    Word records the cursor position of this edit as the "_GoBack"
#    bookmark, which moves here from its previous location (end of the
#    word "mysql" further down in the document). Since a document can
#    only have a single bookmark of a given name, re-adding "_GoBack" at
#    the new spot automatically removes it from the old spot.
# ---------------------------------------------------------------------

$gleison = $d.Content
$gleison.Find.Execute("Gleison Silva ")
$gleisonStart = $gleison.Start

# Temporarily re-seat "_GoBack" right before the trailing " - " that follows
# "...Freire - 1203037" so that run stays split off once the paragraph is
# touched by the deletion below (otherwise it gets merged back into its
# neighbour run).
$dash = $d.Range($gleisonStart, $d.Content.End)
$dash.Find.Execute(" - ")
$dashPoint = $d.Range($dash.Start, $dash.Start)
$d.Bookmarks.Add("_GoBack", $dashPoint)

# Delete "Silva " leaving "Gleison " (with trailing space) as its own run
$silva = $d.Range($gleisonStart, $d.Content.End)
$silva.Find.Execute("Silva ")
$silva.Delete()

# Re-seat the "_GoBack" bookmark at its real, final location immediately
# before "Freire".
$freire = $d.Range($gleisonStart, $d.Content.End)
$freire.Find.Execute("Freire")
$bmPoint = $d.Range($freire.Start, $freire.Start)
$d.Bookmarks.Add("_GoBack", $bmPoint)

# ---------------------------------------------------------------------
# 2) "Taxa Adm Geral (valor)(perda de documento..." run split shifts:
#    "valor" moves from the end of the "Geral (valor" run to the start
#    of the ")(" run, i.e. " Geral (" + "valor)(" instead of
#    " Geral (valor" + ")(".
# ---------------------------------------------------------------------

$parte1 = $d.Content
$parte1.Find.Execute(" Geral (valor")
$parte1.Text = " Geral ("

$parte2 = $d.Content
$parte2.Find.Execute(")(")
$parte2.Text = "valor)("
